$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 8332.666999999999
$ws.Range("I39").Value = 5000
$ws.Range("J39").Value = 9999
$ws.Range("K39").Value = 15000
$ws.Range("L39").Value = 29997
$ws.Range("M39").Value = -14704
$ws.Range("N39").Value = -30589
$ws.Range("H62").Value = 24045.938
$ws.Range("I62").Value = 11340.546
$ws.Range("J62").Value = 51997.8
$ws.Range("K62").Value = 11340.546
$ws.Range("L62").Value = 51997.8
$ws.Range("M62").Value = -10716.546
$ws.Range("N62").Value = -53245.8
$ws.Range("H65").Value = 24045.938
$ws.Range("I65").Value = 11340.546
$ws.Range("J65").Value = 51997.8
$ws.Range("K65").Value = 56702.73
$ws.Range("L65").Value = 259989
$ws.Range("M65").Value = -53582.73
$ws.Range("N65").Value = -266229
$ws.Range("H86").Value = 44591830
$ws.Range("J86").Value = 8001751.5
$ws.Range("L86").Value = 8001751.5
$ws.Range("N86").Value = -8003997.5
$ws.Range("H89").Value = 44591830
$ws.Range("J89").Value = 8001751.5
$ws.Range("L89").Value = 40008757.5
$ws.Range("N89").Value = -40019989.5
$ws.Range("H106").Value = 125003224
$ws.Range("I106").Value = 142860110
$ws.Range("K106").Value = 142860110
$ws.Range("M106").Value = -142859479
$ws.Range("H132").Value = 2471.2917
$ws.Range("I132").Value = 2536.238
$ws.Range("K132").Value = 7608.714
$ws.Range("M132").Value = -5078.714
$ws.Range("H138").Value = 2637440.5
$ws.Range("I138").Value = 2631.3076
$ws.Range("J138").Value = 4007541.2
$ws.Range("K138").Value = 7893.9228
$ws.Range("L138").Value = 12022623.6
$ws.Range("M138").Value = -2753.9228
$ws.Range("N138").Value = -12032903.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 83336640
$ws.Range("I2").Value = 3333.1667
$ws.Range("J2").Value = 166669950
$ws.Range("K2").Value = 3333.1667
$ws.Range("L2").Value = 166669950
$ws.Range("M2").Value = -3220.1667
$ws.Range("N2").Value = -166670176
$ws.Range("H45").Value = 3982.4285
$ws.Range("J45").Value = 3012.75
$ws.Range("L45").Value = 3012.75
$ws.Range("N45").Value = -3766.75
$ws.Range("H61").Value = 32263584
$ws.Range("I61").Value = 1913.6818
$ws.Range("K61").Value = 1913.6818
$ws.Range("M61").Value = -1701.6818
$ws.Range("H116").Value = 83336640
$ws.Range("I116").Value = 3333.1667
$ws.Range("J116").Value = 166669950
$ws.Range("K116").Value = 3333.1667
$ws.Range("L116").Value = 166669950
$ws.Range("M116").Value = -1039.1667
$ws.Range("N116").Value = -166674538
$ws.Range("H122").Value = 5119.646
$ws.Range("I122").Value = 4335.5
$ws.Range("K122").Value = 13006.5
$ws.Range("M122").Value = -10556.5
$ws.Range("H132").Value = 7735.171
$ws.Range("I132").Value = 5935.9165
$ws.Range("K132").Value = 17807.7495
$ws.Range("M132").Value = -15277.7495
$ws.Range("H136").Value = 32263584
$ws.Range("I136").Value = 1913.6818
$ws.Range("K136").Value = 5741.0454
$ws.Range("M136").Value = -3191.0454

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 83336640
$ws.Range("I3").Value = 3333.1667
$ws.Range("J3").Value = 166669950
$ws.Range("K3").Value = 3333.1667
$ws.Range("L3").Value = 166669950
$ws.Range("M3").Value = -3219.1667
$ws.Range("N3").Value = -166670178

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9086.799999999999
$ws.Range("I31").Value = 3673.1428
$ws.Range("K31").Value = 3673.1428
$ws.Range("M31").Value = -3378.1428
$ws.Range("H34").Value = 9086.799999999999
$ws.Range("I34").Value = 3673.1428
$ws.Range("K34").Value = 3673.1428
$ws.Range("M34").Value = -3471.1428
$ws.Range("H99").Value = 7349.625
$ws.Range("I99").Value = 6012
$ws.Range("K99").Value = 6012
$ws.Range("M99").Value = -4514
$ws.Range("H107").Value = 1703.5625
$ws.Range("I107").Value = 398.36365
$ws.Range("J107").Value = 2387.238
$ws.Range("K107").Value = 398.36365
$ws.Range("L107").Value = 2387.238
$ws.Range("M107").Value = 1521.63635
$ws.Range("N107").Value = -6227.237999999999
$ws.Range("H126").Value = 7349.625
$ws.Range("I126").Value = 6012
$ws.Range("K126").Value = 18036
$ws.Range("M126").Value = -15566

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 387136.8
$ws.Range("I70").Value = 730178.2
$ws.Range("J70").Value = 9791.299999999999
$ws.Range("K70").Value = 730178.2
$ws.Range("L70").Value = 9791.299999999999
$ws.Range("M70").Value = -729908.2
$ws.Range("N70").Value = -10331.3
$ws.Range("H73").Value = 387136.8
$ws.Range("I73").Value = 730178.2
$ws.Range("J73").Value = 9791.299999999999
$ws.Range("K73").Value = 730178.2
$ws.Range("L73").Value = 9791.299999999999
$ws.Range("M73").Value = -729242.2
$ws.Range("N73").Value = -11663.3
$ws.Range("H102").Value = 2020.0952
$ws.Range("I102").Value = 1848.4375
$ws.Range("K102").Value = 1848.4375
$ws.Range("M102").Value = -226.4375
$ws.Range("H122").Value = 2860639.5
$ws.Range("I122").Value = 4466180.5
$ws.Range("K122").Value = 13398541.5
$ws.Range("M122").Value = -13396091.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5481.6665
$ws.Range("I7").Value = 4166.8335
$ws.Range("K7").Value = 4166.8335
$ws.Range("M7").Value = -4054.8335
$ws.Range("H40").Value = 7376.923
$ws.Range("I40").Value = 6779.6
$ws.Range("K40").Value = 6779.6
$ws.Range("M40").Value = -6643.6
$ws.Range("H68").Value = 6595.4
$ws.Range("I68").Value = 5996.75
$ws.Range("K68").Value = 5996.75
$ws.Range("M68").Value = -5247.75
$ws.Range("H71").Value = 6595.4
$ws.Range("I71").Value = 5996.75
$ws.Range("K71").Value = 29983.75
$ws.Range("M71").Value = -26239.75
$ws.Range("H122").Value = 3994.4707
$ws.Range("I122").Value = 3303.862
$ws.Range("K122").Value = 9911.585999999999
$ws.Range("M122").Value = -7461.585999999999
$ws.Range("H126").Value = 5481.6665
$ws.Range("I126").Value = 4166.8335
$ws.Range("K126").Value = 12500.5005
$ws.Range("M126").Value = -10030.5005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 8782673
$ws.Range("I81").Value = 694357.2
$ws.Range("J81").Value = 22263200
$ws.Range("K81").Value = 1388714.4
$ws.Range("L81").Value = 44526400
$ws.Range("M81").Value = -1387653.4
$ws.Range("N81").Value = -44528522
$ws.Range("H84").Value = 8782673
$ws.Range("I84").Value = 694357.2
$ws.Range("J84").Value = 22263200
$ws.Range("K84").Value = 6943572
$ws.Range("L84").Value = 222632000
$ws.Range("M84").Value = -6938268
$ws.Range("N84").Value = -222642608
$ws.Range("I107").Value = 284.125
$ws.Range("J107").Value = 19609268
$ws.Range("K107").Value = 852.375
$ws.Range("L107").Value = 58827804
$ws.Range("M107").Value = 1067.625
$ws.Range("N107").Value = -58831644
$ws.Range("H126").Value = 4351.6
$ws.Range("I126").Value = 3762.5557
$ws.Range("J126").Value = 5235.1665
$ws.Range("K126").Value = 11287.6671
$ws.Range("L126").Value = 15705.4995
$ws.Range("M126").Value = -8817.667099999999
$ws.Range("N126").Value = -20645.4995
$ws.Range("H135").Value = 41000
$ws.Range("J135").Value = 41000
$ws.Range("L135").Value = 41000
$ws.Range("N135").Value = -51140
$ws.Range("H136").Value = 30610030
$ws.Range("I136").Value = 52632724
$ws.Range("J136").Value = 722089.9399999999
$ws.Range("K136").Value = 157898172
$ws.Range("L136").Value = 2166269.82
$ws.Range("M136").Value = -157895622
$ws.Range("N136").Value = -2171369.82
